$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("K")

$ws.Rows.Item(3).Insert()
